$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add "Example 5" table (mirrors the structure of the "Example 4" table just
# above it: a merged/rotated label cell in column A, a bold header row, and
# data rows). Re-use the existing "Example 4" header + first five data rows'
# formatting as a template so the new block matches the sheet's look & feel.
# ---------------------------------------------------------------------------
$ws.Range("A20:A25").Merge()
$ws.Range("A10:G15").Copy()
$ws.Range("A20:G25").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 20 - header for Example 5
$ws.Range("A20").Value = "Example 5"
$ws.Range("B20").Value = "Line"
$ws.Range("C20").Value = "Action"
$ws.Range("D20").Value = "Object"
$ws.Range("E20").Value = "Side Effect"
$ws.Range("F20").Value = "Return Value"
$ws.Range("G20").Value = "Is Return Value Used?"

# Row 21
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = "Method call (map)"
$ws.Range("D21").Value = "The outer array"
$ws.Range("E21").Value = "None"
$ws.Range("F21").Value = "New transformed array"
$ws.Range("G21").Value = "No"

# Row 22
$ws.Range("B22").Value = "1 thru 5"
$ws.Range("C22").Value = "Outer block execution"
$ws.Range("D22").Value = "Each sub array"
$ws.Range("E22").Value = "None"
$ws.Range("F22").Value = "New transformed array"
$ws.Range("G22").Value = "Yes, used by original map for transformation"

# Row 23
$ws.Range("B23").Value = 2
$ws.Range("C23").Value = "Method call (map)"
$ws.Range("D23").Value = "Each sub array"
$ws.Range("E23").Value = "None"
$ws.Range("F23").Value = "New transformed array"
$ws.Range("G23").Value = "Yes, used to determine return value of outer block"

# Row 24
$ws.Range("B24").Value = "2 thru 4"
$ws.Range("C24").Value = "Inner block execution"
$ws.Range("D24").Value = "Element of sub array in that iteration"
$ws.Range("E24").Value = "None"
$ws.Range("F24").Value = "Transformed element"
$ws.Range("G24").Value = "Yes, used by inner map for transformation"

# Row 25
$ws.Range("B25").Value = 3
$ws.Range("C25").Value = "num * 2"
$ws.Range("D25").Value = "n/a"
$ws.Range("E25").Value = "None"
$ws.Range("F25").Value = "An integer"
$ws.Range("G25").Value = "Yes, used to determine return value of inner block"

# ---------------------------------------------------------------------------
# Three trailing blank rows (21 -> 28) under the new table, formatted like
# the rotated label column but left-aligned back to general/center so they
# read as blank spacer cells.
# ---------------------------------------------------------------------------
$ws.Range("A4").Copy()
$ws.Range("A26:A28").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A26:A28").Value = ""
$ws.Range("A26:A28").HorizontalAlignment = 1  # xlHAlignGeneral

# ---------------------------------------------------------------------------
# Update the view so the window is scrolled/selected the way it was left
# after the edit (top-left anchored at row 13, selection on the newly typed
# E22 cell).
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("E22").Select()

Write-Host "Added Example 5 notes table (rows 20-28)"
